$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated loading_percent results for the 380 kV case (Case_0_166)
# Rows 2-25 correspond to time steps; columns C-O are the line loading series.
$data = @{
    2 = @{ "C" = 2.267666833824733; "D" = 3.116160134795282; "E" = 8.90790172176087; "F" = 22.84856934594624; "G" = 3.579128696932475; "I" = 18.61228321217679; "M" = 20.10933778132328; "N" = 17.14847856116813; "O" = 19.69620884232447 }
    3 = @{ "C" = 2.236489039930447; "D" = 3.128104637768304; "E" = 9.071111396561825; "F" = 22.40334625795055; "G" = 3.582256344044262; "I" = 18.36336929813508; "M" = 19.3113018106483; "N" = 16.86202298659484; "O" = 19.43208596770203 }
    4 = @{ "C" = 2.217516696863296; "D" = 3.135715199765265; "E" = 9.17603757175806; "F" = 22.13408135302123; "G" = 3.584278229432194; "I" = 18.21543289074031; "M" = 18.80535993785947; "N" = 16.68567178769943; "O" = 19.27516094715822 }
    5 = @{ "C" = 2.209837287193575; "D" = 3.138886508346932; "E" = 9.219986010606579; "F" = 22.02555685186682; "G" = 3.585127778400224; "I" = 18.15645614037657; "M" = 18.59546109400052; "N" = 16.61377914302983; "O" = 19.21261367370104 }
    6 = @{ "C" = 2.208565527891166; "D" = 3.139417337563406; "E" = 9.227355633920139; "F" = 22.00761409548249; "G" = 3.585270395097859; "I" = 18.14674431467211; "M" = 18.56039169043028; "N" = 16.60184236284296; "O" = 19.2023146719334 }
    7 = @{ "C" = 2.217412907361846; "D" = 3.135757685450338; "E" = 9.17662545085725; "F" = 22.13261265662346; "G" = 3.584289582914316; "I" = 18.21463211219791; "M" = 18.80254384188861; "N" = 16.68470221319256; "O" = 19.27431163544018 }
    8 = @{ "C" = 2.256886652045655; "D" = 3.120221430910871; "E" = 8.963200644116535; "F" = 22.69430743640546; "G" = 3.580186099282567; "I" = 18.52548826413983; "M" = 19.83764547887423; "N" = 17.04985843081084; "O" = 19.60410000765757 }
    9 = @{ "C" = 2.335255741407644; "D" = 3.091931299083549; "E" = 8.581873619711793; "F" = 23.82085353474157; "G" = 3.572940329287471; "I" = 19.17042021024251; "M" = 21.7300622148706; "N" = 17.75840759113683; "O" = 20.28872559365047 }
    10 = @{ "C" = 2.392894322853097; "D" = 3.072447446575471; "E" = 8.324092448126947; "F" = 24.65393969638901; "G" = 3.568099410445844; "I" = 19.66095779226536; "M" = 23.02418899172949; "N" = 18.26944999930586; "O" = 20.80970170532845 }
    11 = @{ "C" = 2.419020096673616; "D" = 3.063860622941303; "E" = 8.211617299656963; "F" = 25.03210792775555; "G" = 3.566000683661486; "I" = 19.88670007191656; "M" = 23.58994099738629; "N" = 18.49888101375172; "O" = 21.04950481655393 }
    12 = @{ "C" = 2.428890890410405; "D" = 3.060648344609915; "E" = 8.169710008762834; "F" = 25.17503660226199; "G" = 3.56522072809927; "I" = 19.97246919196362; "M" = 23.8007339437164; "N" = 18.58524528064757; "O" = 21.14062388877539 }
    13 = @{ "C" = 2.426766176569334; "D" = 3.06133842167437; "E" = 8.178705113668663; "F" = 25.14426905863891; "G" = 3.56538804921528; "I" = 19.95398596376138; "M" = 23.75549121083012; "N" = 18.56666935869195; "O" = 21.12098740710349 }
    14 = @{ "C" = 2.419832673171377; "D" = 3.063595560175298; "E" = 8.208155865592495; "F" = 25.04387316085324; "G" = 3.565936220394605; "I" = 19.89375112590419; "M" = 23.60735285755998; "N" = 18.50599703610612; "O" = 21.0569955354768 }
    15 = @{ "C" = 2.415582516280257; "D" = 3.064983237580136; "E" = 8.226284330457169; "F" = 24.98233730882644; "G" = 3.566273914175112; "I" = 19.85689014050916; "M" = 23.51616125962559; "N" = 18.4687639917083; "O" = 21.01783641017624 }
    16 = @{ "C" = 2.391184276618969; "D" = 3.073014145660095; "E" = 8.331538997853183; "F" = 24.62919547431357; "G" = 3.568238641069593; "I" = 19.64624968659666; "M" = 22.98674106415534; "N" = 18.2543880785051; "O" = 20.7940785403046 }
    17 = @{ "C" = 2.376186168787622; "D" = 3.078011381728062; "E" = 8.397333290293199; "F" = 24.41223190146688; "G" = 3.569470366899079; "I" = 19.51763031565045; "M" = 22.65597055343066; "N" = 18.12203850212096; "O" = 20.65746307731851 }
    18 = @{ "C" = 2.367551294271564; "D" = 3.080911706678383; "E" = 8.435627614898458; "F" = 24.28737332233715; "G" = 3.570188563629171; "I" = 19.44390122697095; "M" = 22.46356688300708; "N" = 18.04563238881919; "O" = 20.57915532921727 }
    19 = @{ "C" = 2.36462651426264; "D" = 3.081898192380679; "E" = 8.448671038867991; "F" = 24.24509185078261; "G" = 3.570433408305677; "I" = 19.41898334447281; "M" = 22.39805724048565; "N" = 18.01971658932601; "O" = 20.55269087366669 }
    20 = @{ "C" = 2.377783673541762; "D" = 3.077476724534729; "E" = 8.390282710923175; "F" = 24.43533614460781; "G" = 3.569338240097222; "I" = 19.53129689732207; "M" = 22.69140560279838; "N" = 18.13615712153003; "O" = 20.67197877606482 }
    21 = @{ "C" = 2.421869890619367; "D" = 3.062931518641959; "E" = 8.199486919916597; "F" = 25.07337057599626; "G" = 3.565774808611597; "I" = 19.91143650329869; "M" = 23.65095922886917; "N" = 18.52383258501156; "O" = 21.07578379744986 }
    22 = @{ "C" = 2.450547884358581; "D" = 3.053654655394004; "E" = 8.078779285627764; "F" = 25.48868932317644; "G" = 3.563532048618784; "I" = 20.16150379728794; "M" = 24.25794976634317; "N" = 18.77415810834763; "O" = 21.34146353234912 }
    23 = @{ "C" = 2.435257111343106; "D" = 3.058585043338801; "E" = 8.142839690478077; "F" = 25.26722962439144; "G" = 3.564721197739841; "I" = 20.02791813932818; "M" = 23.93587213268806; "N" = 18.64085795863757; "O" = 21.19953364764385 }
    24 = @{ "C" = 2.377061479169553; "D" = 3.077718357821805; "E" = 8.393468818699432; "F" = 24.42489108816609; "G" = 3.569397943277856; "I" = 19.52511755652527; "M" = 22.67539238069711; "N" = 18.12977507375716; "O" = 20.66541549170612 }
    25 = @{ "C" = 2.314006589404088; "D" = 3.099354204761634; "E" = 8.681081193045891; "F" = 23.51444389918862; "G" = 3.574815340379839; "I" = 18.99266077954127; "M" = 21.23422152566114; "N" = 17.56805320244677; "O" = 20.09998255463013 }
}

foreach ($rowNum in $data.Keys) {
    $rowValues = $data[$rowNum]
    foreach ($col in $rowValues.Keys) {
        $ws.Range("$col$rowNum").Value = $rowValues[$col]
    }
}

Write-Host "Updated $($data.Count) rows of loading_percent data"
